$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.084.01'
$ws.Range("E2").Value = '  -0.33%  '

# Row 3
$ws.Range("D3").Value = '3.098.69'
$ws.Range("E3").Value = '  -0.08%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.65'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -1.15%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.97'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +3.57%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("D8").Value = '3.097.75'
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.513'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -1.00%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.38'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -1.30%  '

# Row 11
$ws.Range("E11").Value = '  -0.62%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -1.99%  '

# Row 13
$ws.Range("E13").Value = '  -2.00%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.08'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -1.42%  '

# Row 15
$ws.Range("E15").Value = '  +0.35%  '

# Row 16
$ws.Range("D16").Value = '3.617.52'
$ws.Range("E16").Value = '  +0.05%  '

# Row 17
$ws.Range("D17").Value = '67.039.04'
$ws.Range("E17").Value = '  -0.28%  '

# Row 18
$ws.Range("E18").Value = '  -1.04%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.73'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +0.00%  '

# Row 20
$ws.Range("D20").Value = '3.102.23'
$ws.Range("E20").Value = '  +0.22%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '489.02'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +0.03%  '

# Row 22
$ws.Range("E22").Value = '  -1.08%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.687'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -1.40%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.46'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -0.21%  '

# Row 25
$ws.Range("E25").Value = '  -0.14%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.58'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -3.39%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.16'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -2.69%  '

# Row 28
$ws.Range("E28").Value = '  -0.06%  '

# Row 29
$ws.Range("E29").Value = '  +1.48%  '

# Row 30
$ws.Range("E30").Value = '  -1.41%  '

# Row 31
$ws.Range("E31").Value = '  -2.36%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.13'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -0.39%  '

# Row 33
$ws.Range("E33").Value = '  -1.04%  '

# Row 34
$ws.Range("E34").Value = '  -0.05%  '

# Row 35
$ws.Range("E35").Value = '  +0.03%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '47.58'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +2.59%  '

# Row 37
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.57'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -3.43%  '

# Row 38
$ws.Range("B38").Value = 'Mantle'
$ws.Range("C38").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.945'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -2.27%  '

# Row 39
$ws.Range("E39").Value = '  +2.29%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.03'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +0.63%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '49.13'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -1.31%  '

# Row 42
$ws.Range("E42").Value = '  -0.35%  '

# Row 43
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.72'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +5.74%  '

# Row 44
$ws.Range("B44").Value = 'Cosmos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.24'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -2.28%  '

# Row 45
$ws.Range("D45").Value = '2.802.15'
$ws.Range("E45").Value = '  +0.32%  '

# Row 46
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '369.53'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -2.84%  '

# Row 47
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0345'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -1.31%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.68'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -0.01%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.58'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +2.92%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.29'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +4.91%  '
